$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new row for BERT (becomes new row 2; pushes Donut and below down by one) ---
$ws.Rows.Item(2).Insert()

# --- Insert new row for BART (becomes new row 4, between Donut and Attention) ---
$ws.Rows.Item(4).Insert()

# Fill row 2 (BERT)
$ws.Cells.Item(2,1).Value = 'BERT: Pre-training of Deep Bidirectional Transformers for Language Understanding'
$ws.Cells.Item(2,2).Value = 'May'
$ws.Cells.Item(2,3).Value = '2019'
$ws.Cells.Item(2,4).Value = 'Devlin, Jacob and Chang, Ming-Wei and Lee, Kenton and Toutanova, Kristina'
$ws.Cells.Item(2,5).Value = 'arXiv'
$ws.Cells.Item(2,6).Value = 'Devlin et al. - 2019 - BERT Pre-training of Deep Bidirectional Transform.pdf'
$ws.Cells.Item(2,7).Value = 'Computer Science - Computation and Language'
$ws.Cells.Item(2,9).Value = 'We introduce a new language representation model called BERT, which stands for Bidirectional Encoder Representations from Transformers. Unlike recent language representation models, BERT is designed to pre-train deep bidirectional representations from unlabeled text by jointly conditioning on both left and right context in all layers. As a result, the pre-trained BERT model can be fine-tuned with just one additional output layer to create state-of-the-art models for a wide range of tasks, such as question answering and language inference, without substantial task-specific architecture modifications. BERT is conceptually simple and empirically powerful. It obtains new state-of-the-art results on eleven natural language processing tasks, including pushing the GLUE score to 80.5% (7.7% point absolute improvement), MultiNLI accuracy to 86.7% (4.6% absolute improvement), SQuAD v1.1 question answering Test F1 to 93.2 (1.5 point absolute improvement) and SQuAD v2.0 Test F1 to 83.1 (5.1 point absolute improvement).'

# Fill row 4 (BART)
$ws.Cells.Item(4,1).Value = 'BART: Denoising Sequence-to-Sequence Pre-training for Natural Language Generation, Translation, and Comprehension'
$ws.Cells.Item(4,2).Value = 'October'
$ws.Cells.Item(4,3).Value = '2019'
$ws.Cells.Item(4,4).Value = 'Lewis, Mike and Liu, Yinhan and Goyal, Naman and Ghazvininejad, Marjan and Mohamed, Abdelrahman and Levy, Omer and Stoyanov, Ves and Zettlemoyer, Luke'
$ws.Cells.Item(4,5).Value = 'arXiv'
$ws.Cells.Item(4,6).Value = 'Lewis et al. - 2019 - BART Denoising Sequence-to-Sequence Pre-training .pdf'
$ws.Cells.Item(4,7).Value = 'Computer Science - Computation and Language,Computer Science - Machine Learning,Statistics - Machine Learning'
$ws.Cells.Item(4,9).Value = 'We present BART, a denoising autoencoder for pretraining sequence-to-sequence models. BART is trained by (1) corrupting text with an arbitrary noising function, and (2) learning a model to reconstruct the original text. It uses a standard Tranformer-based neural machine translation architecture which, despite its simplicity, can be seen as generalizing BERT (due to the bidirectional encoder), GPT (with the left-to-right decoder), and many other more recent pretraining schemes. We evaluate a number of noising approaches, finding the best performance by both randomly shuffling the order of the original sentences and using a novel in-filling scheme, where spans of text are replaced with a single mask token. BART is particularly effective when fine tuned for text generation but also works well for comprehension tasks. It matches the performance of RoBERTa with comparable training resources on GLUE and SQuAD, achieves new state-of-the-art results on a range of abstractive dialogue, question answering, and summarization tasks, with gains of up to 6 ROUGE. BART also provides a 1.1 BLEU increase over a back-translation system for machine translation, with only target language pretraining. We also report ablation experiments that replicate other pretraining schemes within the BART framework, to better measure which factors most influence end-task performance.'

# --- Fix formatting for newly-inserted rows (2 and 4) to match existing data rows ---
foreach ($r in @(2,4)) {
  $a = $ws.Cells.Item($r,1)
  $a.Font.Bold = $true
  $a.HorizontalAlignment = -4108
  $a.VerticalAlignment = -4108
  $a.WrapText = $true
  $a.Borders.Item(10).LineStyle = 1
  for ($col = 2; $col -le 9; $col++) {
    $c = $ws.Cells.Item($r, $col)
    $c.Font.Bold = $false
    $c.HorizontalAlignment = -4131
    $c.VerticalAlignment = -4108
    $c.WrapText = $true
    $c.IndentLevel = 1
    $c.Borders.Item(10).LineStyle = 1
  }
  $ws.Rows.Item($r).RowHeight = 96
}

# --- Adjust column A width to match new longest title ---
$ws.Columns.Item(1).ColumnWidth = 29.85546875
